$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.038.18'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '3.797.63'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.18'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.90'
$ws.Range("E6").Value = '  -3.63%  '

$ws.Range("D7").Value = '3.797.54'
$ws.Range("E7").Value = '  +1.30%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E10").Value = '  +1.63%  '

$ws.Range("E11").Value = '  -0.75%  '

$ws.Range("E12").Value = '  -1.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.20'
$ws.Range("E13").Value = '  -3.06%  '

$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").Value = '4.432.49'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").Value = '3.776.38'
$ws.Range("E16").Value = '  +0.78%  '

$ws.Range("D17").Value = '69.140.47'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.38'
$ws.Range("E18").Value = '  +1.24%  '

$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.30'
$ws.Range("E20").Value = '  +1.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.23'
$ws.Range("E21").Value = '  +3.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.95'
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.720'
$ws.Range("E23").Value = '  -1.12%  '

$ws.Range("E24").Value = '  -3.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.56'
$ws.Range("E25").Value = '  -0.95%  '

$ws.Range("E26").Value = '  -3.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  -3.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.96'
$ws.Range("E30").Value = '  -0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.04'
$ws.Range("E31").Value = '  +0.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.39'
$ws.Range("E32").Value = '  -5.84%  '

$ws.Range("D33").Value = '3.947.40'
$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("E34").Value = '  -0.34%  '

$ws.Range("D35").Value = '3.745.45'
$ws.Range("E35").Value = '  +1.67%  '

$ws.Range("E36").Value = '  -2.14%  '

$ws.Range("E37").Value = '  +5.24%  '

$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.89'
$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.322'
$ws.Range("E41").Value = '  -0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.47'
$ws.Range("E43").Value = '  -0.75%  '

$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '418.20'
$ws.Range("E45").Value = '  -4.59%  '

$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.88'
$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.820.25'
$ws.Range("E49").Value = '  +1.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.50'
$ws.Range("E50").Value = '  -2.57%  '
